# Data Drive Invalid Login Test Script
# - Duplicate the "ValidLogin" sheet into a new "InvalidLogin" sheet
# - Replace the duplicated data row with the invalid-login data (abc / xyz)
# - Draw a thin box border around the A1:B2 data grid on both sheets
# - Widen column A to fit the header text on both sheets
# - Make "InvalidLogin" the active/selected sheet, with cell B3 selected (just
#   below the data), and reset both sheets' zoom back to normal (100%)

$wb = $excel.ActiveWorkbook

$validSheet = $wb.Worksheets.Item(1)

# Duplicate the ValidLogin worksheet right after itself, then rename the copy.
$validSheet.Copy($null, $validSheet)
$invalidSheet = $wb.Worksheets.Item(2)
$invalidSheet.Name = "InvalidLogin"

# Give the new sheet its own (invalid) credentials.
$invalidSheet.Range("A2").Value = "abc"
$invalidSheet.Range("B2").Value = "xyz"

# Box the 2x2 data range with a thin border on both sheets.
$validSheet.Range("A1:B2").Borders.LineStyle = 1
$validSheet.Range("A1:B2").Borders.Weight = 2

$invalidSheet.Range("A1:B2").Borders.LineStyle = 1
$invalidSheet.Range("A1:B2").Borders.Weight = 2

# Fit column A to the (now visible) header/username text on both sheets.
$validSheet.Columns.Item(1).ColumnWidth = 9.43
$invalidSheet.Columns.Item(1).ColumnWidth = 9.43

# Reset zoom + selection on the original sheet.
$validSheet.Activate()
$excel.ActiveWindow.Zoom = 100
$validSheet.Range("A1:B2").Select()

# Make InvalidLogin the active/visible tab, reset its zoom, select B3.
$invalidSheet.Activate()
$excel.ActiveWindow.Zoom = 100
$invalidSheet.Range("B3").Select()
